$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H103" = 1424.75
    "I103" = 999
    "J103" = 1566.6666
    "K103" = 2997
    "L103" = 4699.9998
    "M103" = -2411
    "N103" = -5871.9998
    "H111" = 7497.3335
    "I111" = 7250.125
    "J111" = 9475
    "K111" = 21750.375
    "L111" = 28425
    "M111" = -18683.375
    "N111" = -34559
    "H116" = 3326.111
    "I116" = 2727.0667
    "J116" = 4074.9167
    "K116" = 2727.0667
    "L116" = 4074.9167
    "M116" = 714.9333000000001
    "N116" = -10958.9167
    "H132" = 2162.8333
    "I132" = 2011.625
    "K132" = 6034.875
    "M132" = -3504.875
    "H137" = 1070.7142
    "I137" = 0
    "K137" = 0
    "H141" = 2286.75
    "I141" = 1457
    "K141" = 4371
    "M141" = 809
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in @("M137")) {
    $ws.Range($addr).ClearContents()
}

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H38" = 100798.8
    "I38" = 113750
    "J38" = 48994
    "K38" = 113750
    "L38" = 48994
    "M38" = -113283
    "N38" = -49928
    "H110" = 1515.8667
    "I110" = 1619.8334
    "J110" = 1100
    "K110" = 1619.8334
    "L110" = 1100
    "M110" = 425.1666
    "N110" = -5190
    "H122" = 35365.89
    "I122" = 35365.89
    "K122" = 106097.67
    "M122" = -103647.67
    "H132" = 1731.2941
    "I132" = 1631.2142
    "J132" = 2198.3333
    "K132" = 4893.642599999999
    "L132" = 6594.999899999999
    "M132" = -2363.642599999999
    "N132" = -11654.9999
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H86" = 7701.4287
    "I86" = 7701.4287
    "K86" = 7701.4287
    "M86" = -6578.4287
    "H89" = 7701.4287
    "I89" = 7701.4287
    "K89" = 38507.14350000001
    "M89" = -32891.14350000001
    "H99" = 3949.75
    "I99" = 4899.5
    "K99" = 4899.5
    "M99" = -3401.5
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H22" = 1934.6666
    "I22" = 0
    "K22" = 0
    "H31" = 1052.8
    "I31" = 860
    "K31" = 860
    "M31" = -565
    "H34" = 1052.8
    "I34" = 860
    "K34" = 860
    "M34" = -658
    "H58" = 2505.889
    "I58" = 1847
    "K58" = 1847
    "M58" = -1644
    "H86" = 4794.9
    "I86" = 4678.7144
    "K86" = 4678.7144
    "M86" = -3555.7144
    "H89" = 4794.9
    "I89" = 4678.7144
    "K89" = 23393.572
    "M89" = -17777.572
    "H132" = 7206.6875
    "I132" = 7206.6875
    "K132" = 21620.0625
    "M132" = -19090.0625
    "H134" = 2202.182
    "I134" = 2232.4
    "K134" = 6697.200000000001
    "M134" = -4162.200000000001
    "H136" = 2505.889
    "I136" = 1847
    "K136" = 5541
    "M136" = -2991
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in @("M22")) {
    $ws.Range($addr).ClearContents()
}

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H38" = 1796.1
    "I38" = 2401.5715
    "J38" = 383.33334
    "K38" = 7204.7145
    "L38" = 1150.00002
    "M38" = -6857.7145
    "N38" = -1844.00002
    "H75" = 2437.6667
    "I75" = 2449
    "K75" = 7347
    "M75" = -6349
    "H78" = 2437.6667
    "I78" = 2449
    "K78" = 22041
    "M78" = -17049
    "H86" = 199
    "J86" = 197
    "L86" = 591
    "N86" = -2963
    "H89" = 199
    "J89" = 197
    "L89" = 1773
    "N89" = -13629
    "H131" = 401838.38
    "I131" = 1269
    "J131" = 528333.9399999999
    "K131" = 3807
    "L131" = 1585001.82
    "M131" = 1233
    "N131" = -1595081.82
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H70" = 0
    "J70" = 0
    "N70" = 0
    "H73" = 0
    "J73" = 0
    "N73" = 0
    "H80" = 6820.2
    "I80" = 4523.75
    "K80" = 4523.75
    "M80" = -3525.75
    "H83" = 6820.2
    "I83" = 4523.75
    "K83" = 22618.75
    "M83" = -17626.75
    "H93" = 0
    "J93" = 0
    "N93" = 0
    "H98" = 13995.4
    "J98" = 13995.4
    "L98" = 13995.4
    "N98" = -19985.4
    "H122" = 2799.25
    "I122" = 3699
    "K122" = 11097
    "M122" = -8647
    "H126" = 5297.222
    "I126" = 3320.6667
    "J126" = 6285.5
    "K126" = 9962.000100000001
    "L126" = 18856.5
    "M126" = -7492.000100000001
    "N126" = -23796.5
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in @("L70", "L73", "L93")) {
    $ws.Range($addr).ClearContents()
}

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H16" = 0
    "I16" = 0
    "K16" = 0
    "H40" = 4024.5
    "I40" = 0
    "K40" = 0
    "H46" = 2267.1
    "I46" = 2038.8572
    "J46" = 2799.6667
    "K46" = 2038.8572
    "L46" = 2799.6667
    "M46" = -1850.8572
    "N46" = -3175.6667
    "H100" = 5416.5
    "I100" = 5416.5
    "K100" = 5416.5
    "M100" = -4875.5
    "H136" = 1498.8572
    "I136" = 1498.8572
    "K136" = 4496.571599999999
    "M136" = -1946.571599999999
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in @("M16", "M40")) {
    $ws.Range($addr).ClearContents()
}

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H114" = 30000
    "J114" = 30000
    "L114" = 30000
    "N114" = -38678
    "H126" = 2000
    "J126" = 2000
    "L126" = 6000
    "N126" = -10940
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
